# Adds EXPERIMENT 3 / HYPERPARAMETER SEARCH 5 block (rows 126-162) and a fresh
# "CHOSEN ARCHITECTURE" section (rows 164-165) to Sheet1, mirroring the formatting
# of the existing analogous sections elsewhere on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Section banner: "EXPERIMENT 3"  (row 126, yellow band like row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Copy()
$ws.Range("A126:I126").PasteSpecial(-4122, 0)
$ws.Range("A126:I126").Merge()
$ws.Range("A126").Value = "EXPERIMENT 3"

# ---------------------------------------------------------------------------
# Section banner: "HYPERPARAMETER SEARCH 5"  (row 127, grey band like row 117)
# ---------------------------------------------------------------------------
$ws.Range("A117:I117").Copy()
$ws.Range("A127:I127").PasteSpecial(-4122, 0)
$ws.Range("A127:I127").Merge()
$ws.Range("A127").Value = "HYPERPARAMETER SEARCH 5"

# ---------------------------------------------------------------------------
# Column header row 128 (same layout as header rows 3 / 118)
# ---------------------------------------------------------------------------
$ws.Range("A118:I118").Copy()
$ws.Range("A128:I128").PasteSpecial(-4122, 0)
$ws.Range("A128").Value = "Alpha"
$ws.Range("B128").Value = "Lambda"
$ws.Range("C128").Value = "Name"
$ws.Range("D128").Value = "Worker"
$ws.Range("E128").Value = "Val Log Loss"
$ws.Range("F128").Value = "Best IOU"
$ws.Range("G128").Value = "Best Threshold"
$ws.Range("H128").Value = "Threshold selection"
$ws.Range("I128").Value = "Notes"

# ---------------------------------------------------------------------------
# Data rows 129-158: 30 hyperparameter-search runs (run85..run114)
# Formatting copied from row 85 (A/B = scientific-notation style, C-I = general)
# ---------------------------------------------------------------------------
$ws.Range("A85:I85").Copy()
$ws.Range("A129:I129").PasteSpecial(-4122, 0)
$ws.Range("A129").Value = 0.0000015654
$ws.Range("B129").Value = 0.26921
$ws.Range("C129").Value = "run85"
$ws.Range("D129").Value = 1

$ws.Range("A85:I85").Copy()
$ws.Range("A130:I130").PasteSpecial(-4122, 0)
$ws.Range("A130").Value = 0.0000035554
$ws.Range("B130").Value = 0.000030387
$ws.Range("C130").Value = "run86"
$ws.Range("D130").Value = 2

$ws.Range("A85:I85").Copy()
$ws.Range("A131:I131").PasteSpecial(-4122, 0)
$ws.Range("A131").Value = 0.000011537
$ws.Range("B131").Value = 0.2741
$ws.Range("C131").Value = "run87"
$ws.Range("D131").Value = 3

$ws.Range("A85:I85").Copy()
$ws.Range("A132:I132").PasteSpecial(-4122, 0)
$ws.Range("A132").Value = 0.000011713
$ws.Range("B132").Value = 0.00073872
$ws.Range("C132").Value = "run88"
$ws.Range("D132").Value = 4

$ws.Range("A85:I85").Copy()
$ws.Range("A133:I133").PasteSpecial(-4122, 0)
$ws.Range("A133").Value = 0.000013078
$ws.Range("B133").Value = 0.000055114
$ws.Range("C133").Value = "run89"
$ws.Range("D133").Value = 5

$ws.Range("A85:I85").Copy()
$ws.Range("A134:I134").PasteSpecial(-4122, 0)
$ws.Range("A134").Value = 0.000016275
$ws.Range("B134").Value = 0.00007581
$ws.Range("C134").Value = "run90"
$ws.Range("D134").Value = 6

$ws.Range("A85:I85").Copy()
$ws.Range("A135:I135").PasteSpecial(-4122, 0)
$ws.Range("A135").Value = 0.000018201
$ws.Range("B135").Value = 0.000050484
$ws.Range("C135").Value = "run91"
$ws.Range("D135").Value = 7

$ws.Range("A85:I85").Copy()
$ws.Range("A136:I136").PasteSpecial(-4122, 0)
$ws.Range("A136").Value = 0.000020148
$ws.Range("B136").Value = 0.0012329
$ws.Range("C136").Value = "run92"
$ws.Range("D136").Value = 8

$ws.Range("A85:I85").Copy()
$ws.Range("A137:I137").PasteSpecial(-4122, 0)
$ws.Range("A137").Value = 0.000022684
$ws.Range("B137").Value = 0.0007946
$ws.Range("C137").Value = "run93"
$ws.Range("D137").Value = 9

$ws.Range("A85:I85").Copy()
$ws.Range("A138:I138").PasteSpecial(-4122, 0)
$ws.Range("A138").Value = 0.00003361
$ws.Range("B138").Value = 0.00052887
$ws.Range("C138").Value = "run94"
$ws.Range("D138").Value = 10

$ws.Range("A85:I85").Copy()
$ws.Range("A139:I139").PasteSpecial(-4122, 0)
$ws.Range("A139").Value = 0.000038006
$ws.Range("B139").Value = 0.00048127
$ws.Range("C139").Value = "run95"
$ws.Range("D139").Value = 11

$ws.Range("A85:I85").Copy()
$ws.Range("A140:I140").PasteSpecial(-4122, 0)
$ws.Range("A140").Value = 0.000040787
$ws.Range("B140").Value = 0.001023
$ws.Range("C140").Value = "run96"
$ws.Range("D140").Value = 12

$ws.Range("A85:I85").Copy()
$ws.Range("A141:I141").PasteSpecial(-4122, 0)
$ws.Range("A141").Value = 0.000041844
$ws.Range("B141").Value = 0.00084802
$ws.Range("C141").Value = "run97"
$ws.Range("D141").Value = 13

$ws.Range("A85:I85").Copy()
$ws.Range("A142:I142").PasteSpecial(-4122, 0)
$ws.Range("A142").Value = 0.000065612
$ws.Range("B142").Value = 0.003607
$ws.Range("C142").Value = "run98"
$ws.Range("D142").Value = 14

$ws.Range("A85:I85").Copy()
$ws.Range("A143:I143").PasteSpecial(-4122, 0)
$ws.Range("A143").Value = 0.00006827
$ws.Range("B143").Value = 0.18602
$ws.Range("C143").Value = "run99"
$ws.Range("D143").Value = 15

$ws.Range("A85:I85").Copy()
$ws.Range("A144:I144").PasteSpecial(-4122, 0)
$ws.Range("A144").Value = 0.00011998
$ws.Range("B144").Value = 0.045356
$ws.Range("C144").Value = "run100"
$ws.Range("D144").Value = 1

$ws.Range("A85:I85").Copy()
$ws.Range("A145:I145").PasteSpecial(-4122, 0)
$ws.Range("A145").Value = 0.000143
$ws.Range("B145").Value = 0.00028769
$ws.Range("C145").Value = "run101"
$ws.Range("D145").Value = 2

$ws.Range("A85:I85").Copy()
$ws.Range("A146:I146").PasteSpecial(-4122, 0)
$ws.Range("A146").Value = 0.00025927
$ws.Range("B146").Value = 0.00014377
$ws.Range("C146").Value = "run102"
$ws.Range("D146").Value = 3

$ws.Range("A85:I85").Copy()
$ws.Range("A147:I147").PasteSpecial(-4122, 0)
$ws.Range("A147").Value = 0.00037919
$ws.Range("B147").Value = 0.0029846
$ws.Range("C147").Value = "run103"
$ws.Range("D147").Value = 4

$ws.Range("A85:I85").Copy()
$ws.Range("A148:I148").PasteSpecial(-4122, 0)
$ws.Range("A148").Value = 0.00039582
$ws.Range("B148").Value = 0.027861
$ws.Range("C148").Value = "run104"
$ws.Range("D148").Value = 5

$ws.Range("A85:I85").Copy()
$ws.Range("A149:I149").PasteSpecial(-4122, 0)
$ws.Range("A149").Value = 0.00055286
$ws.Range("B149").Value = 0.13127
$ws.Range("C149").Value = "run105"
$ws.Range("D149").Value = 6

$ws.Range("A85:I85").Copy()
$ws.Range("A150:I150").PasteSpecial(-4122, 0)
$ws.Range("A150").Value = 0.0016262
$ws.Range("B150").Value = 0.000028017
$ws.Range("C150").Value = "run106"
$ws.Range("D150").Value = 7

$ws.Range("A85:I85").Copy()
$ws.Range("A151:I151").PasteSpecial(-4122, 0)
$ws.Range("A151").Value = 0.0016451
$ws.Range("B151").Value = 0.01021
$ws.Range("C151").Value = "run107"
$ws.Range("D151").Value = 8

$ws.Range("A85:I85").Copy()
$ws.Range("A152:I152").PasteSpecial(-4122, 0)
$ws.Range("A152").Value = 0.0018396
$ws.Range("B152").Value = 0.004497
$ws.Range("C152").Value = "run108"
$ws.Range("D152").Value = 9

$ws.Range("A85:I85").Copy()
$ws.Range("A153:I153").PasteSpecial(-4122, 0)
$ws.Range("A153").Value = 0.0021531
$ws.Range("B153").Value = 0.0044713
$ws.Range("C153").Value = "run109"
$ws.Range("D153").Value = 10

$ws.Range("A85:I85").Copy()
$ws.Range("A154:I154").PasteSpecial(-4122, 0)
$ws.Range("A154").Value = 0.00378
$ws.Range("B154").Value = 0.19318
$ws.Range("C154").Value = "run110"
$ws.Range("D154").Value = 11

$ws.Range("A85:I85").Copy()
$ws.Range("A155:I155").PasteSpecial(-4122, 0)
$ws.Range("A155").Value = 0.0037956
$ws.Range("B155").Value = 0.0014238
$ws.Range("C155").Value = "run111"
$ws.Range("D155").Value = 12

$ws.Range("A85:I85").Copy()
$ws.Range("A156:I156").PasteSpecial(-4122, 0)
$ws.Range("A156").Value = 0.0044606
$ws.Range("B156").Value = 0.11008
$ws.Range("C156").Value = "run112"
$ws.Range("D156").Value = 13

$ws.Range("A85:I85").Copy()
$ws.Range("A157:I157").PasteSpecial(-4122, 0)
$ws.Range("A157").Value = 0.0058469
$ws.Range("B157").Value = 0.000093986
$ws.Range("C157").Value = "run113"
$ws.Range("D157").Value = 14

$ws.Range("A85:I85").Copy()
$ws.Range("A158:I158").PasteSpecial(-4122, 0)
$ws.Range("A158").Value = 0.0094124
$ws.Range("B158").Value = 0.000043753
$ws.Range("C158").Value = "run114"
$ws.Range("D158").Value = 15

# A handful of data rows carry a stray (empty) scientific-format cell in column E,
# mirroring the source template quirk - replicate purely cosmetically.
$ws.Range("A85").Copy()
$ws.Range("E129").PasteSpecial(-4122, 0)
$ws.Range("E131").PasteSpecial(-4122, 0)
$ws.Range("E135").PasteSpecial(-4122, 0)
$ws.Range("E136").PasteSpecial(-4122, 0)
$ws.Range("E145").PasteSpecial(-4122, 0)
$ws.Range("E147").PasteSpecial(-4122, 0)
$ws.Range("E148").PasteSpecial(-4122, 0)
$ws.Range("E151").PasteSpecial(-4122, 0)
$ws.Range("E152").PasteSpecial(-4122, 0)
$ws.Range("E156").PasteSpecial(-4122, 0)
$ws.Range("E157").PasteSpecial(-4122, 0)

# ---------------------------------------------------------------------------
# "Initial range" notes block (rows 160-162), same layout as rows 113-115
# ---------------------------------------------------------------------------
$ws.Range("B160").Value = "Initial range"
$ws.Range("A161").Value = "Alpha:"
$ws.Range("B161").Value = "10 .^ unifrnd(-6,-2,30,1)"
$ws.Range("A162").Value = "Lambda"
$ws.Range("B162").Value = "10 .^ unifrnd(-5, 0, 30, 1)"

# ---------------------------------------------------------------------------
# New "CHOSEN ARCHITECTURE" section (rows 164-165), same layout as rows 117-118
# ---------------------------------------------------------------------------
$ws.Range("A117:I117").Copy()
$ws.Range("A164:I164").PasteSpecial(-4122, 0)
$ws.Range("A164:I164").Merge()
$ws.Range("A164").Value = "CHOSEN ARCHITECTURE"

$ws.Range("A118:I118").Copy()
$ws.Range("A165:I165").PasteSpecial(-4122, 0)
$ws.Range("A165").Value = "Alpha"
$ws.Range("B165").Value = "Lambda"
$ws.Range("C165").Value = "Name"
$ws.Range("D165").Value = "Worker"
$ws.Range("E165").Value = "Val Log Loss"
$ws.Range("F165").Value = "Best IOU"
$ws.Range("G165").Value = "Best Threshold"
$ws.Range("H165").Value = "Threshold selection"
$ws.Range("I165").Value = "Notes"

# ---------------------------------------------------------------------------
# View state: scroll position + active selection, matching the post-edit view
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 110
$win.ScrollColumn = 1
$ws.Range("E155").Select()

